$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.906.59"
$ws.Range("E2").Value = "  +0.10%  "
$ws.Range("D3").Value = "1.638.06"
$ws.Range("E3").Value = "  -0.18%  "
$ws.Range("D4").Value = "'1.002"
$ws.Range("E4").Value = "  -0.41%  "
$ws.Range("D5").Value = "'214.91"
$ws.Range("E5").Value = "  -0.40%  "
$ws.Range("D6").Value = "'0.5069"
$ws.Range("E6").Value = "  +0.80%  "
$ws.Range("D7").Value = "'1.003"
$ws.Range("E7").Value = "  -0.23%  "
$ws.Range("D8").Value = "'0.2558"
$ws.Range("E8").Value = "  -0.43%  "
$ws.Range("D9").Value = "'0.06359"
$ws.Range("E9").Value = "  -0.18%  "
$ws.Range("D10").Value = "'19.45"
$ws.Range("E10").Value = "  -0.75%  "
$ws.Range("D11").Value = "'0.07739"
$ws.Range("E11").Value = "  -0.19%  "
$ws.Range("D12").Value = "'4.276"
$ws.Range("E12").Value = "  +0.59%  "
$ws.Range("D13").Value = "1.636.93"
$ws.Range("E13").Value = "  -0.34%  "
$ws.Range("E14").Value = "  -0.29%  "
$ws.Range("D15").Value = "0.0₅7807"
$ws.Range("E15").Value = "  -0.91%  "
$ws.Range("D16").Value = "'64.25"
$ws.Range("E16").Value = "  +0.32%  "
$ws.Range("D17").Value = "25.947.54"
$ws.Range("E17").Value = "  +0.13%  "
$ws.Range("D18").Value = "'1.002"
$ws.Range("E18").Value = "  -0.30%  "
$ws.Range("D19").Value = "'197.12"
$ws.Range("E19").Value = "  -2.41%  "
$ws.Range("D20").Value = "'4.442"
$ws.Range("E20").Value = "  +1.35%  "
$ws.Range("D21").Value = "'9.929"
$ws.Range("E21").Value = "  +0.42%  "
$ws.Range("D22").Value = "'6.032"
$ws.Range("E22").Value = "  +1.01%  "
$ws.Range("E23").Value = "  -0.24%  "
$ws.Range("D24").Value = "'1.879"
$ws.Range("E24").Value = "  -0.95%  "
$ws.Range("D25").Value = "'141.15"
$ws.Range("E25").Value = "  +0.30%  "
$ws.Range("D26").Value = "'0.1172"
$ws.Range("E26").Value = "  +3.50%  "
$ws.Range("D27").Value = "'6.865"
$ws.Range("E27").Value = "  +1.43%  "
$ws.Range("D28").Value = "'15.68"
$ws.Range("E28").Value = "  +0.44%  "
$ws.Range("D29").Value = "'1.236"
$ws.Range("D30").Value = "'0.04991"
$ws.Range("E30").Value = "  +0.66%  "
$ws.Range("D31").Value = "'3.249"
$ws.Range("E31").Value = "  -0.29%  "
$ws.Range("D32").Value = "'3.183"
$ws.Range("E32").Value = "  -0.24%  "
$ws.Range("D33").Value = "'1.536"
$ws.Range("E33").Value = "  -0.52%  "
$ws.Range("D34").Value = "'2.364"
$ws.Range("E34").Value = "  -0.42%  "
$ws.Range("D35").Value = "'0.8925"
$ws.Range("E35").Value = "  +0.15%  "
$ws.Range("D36").Value = "'2.581"
$ws.Range("E36").Value = "  -1.93%  "
$ws.Range("D37").Value = "1.125.98"
$ws.Range("E37").Value = "  -1.90%  "
$ws.Range("D38").Value = "'0.5440"
$ws.Range("E38").Value = "  -3.33%  "
$ws.Range("D39").Value = "'0.01555"
$ws.Range("E39").Value = "  -0.63%  "
$ws.Range("D40").Value = "'1.002"
$ws.Range("E40").Value = "  -0.39%  "
$ws.Range("B41").Value = "mCoin"
$ws.Range("C41").Value = "https://coinranking.com/coin/fzVgyjBcRc9+mcoin-mcoin"
$ws.Range("D41").Value = "'2.545"
$ws.Range("E41").Value = "  -1.18%  "
$ws.Range("B42").Value = "BabyDogeCoin"
$ws.Range("C42").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D42").Value = "0.0₈129"
$ws.Range("E42").Value = "  +11.13%  "
$ws.Range("D43").Value = "'5.580"
$ws.Range("E43").Value = "  -1.57%  "
$ws.Range("D44").Value = "'0.8162"
$ws.Range("E44").Value = "  +1.27%  "
$ws.Range("E45").Value = "  -0.10%  "
$ws.Range("D46").Value = "1.776.93"
$ws.Range("E46").Value = "  -0.04%  "
$ws.Range("D47").Value = "'0.4531"
$ws.Range("E47").Value = "  -0.41%  "
$ws.Range("D48").Value = "'1.001"
$ws.Range("E48").Value = "  -0.44%  "
$ws.Range("D49").Value = "'54.67"
$ws.Range("E49").Value = "  -0.15%  "
$ws.Range("D50").Value = "'0.05069"
$ws.Range("E50").Value = "  +0.24%  "
$ws.Range("D51").Value = "'1.004"
$ws.Range("E51").Value = "  -0.01%  "
